# Add a new registration row (row 13) below the existing data, matching
# the other rows in columns A:D (Name, Phone Number, Email ID, Event1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(13, 1).Value = "sdfdhjkfgdf"

# Phone numbers in this sheet are stored as text (not numbers), so force
# column B's new cell to text formatting before assigning the value -
# otherwise Excel would auto-detect the numeric-looking string as a number.
$ws.Cells.Item(13, 2).NumberFormat = "@"
$ws.Cells.Item(13, 2).Value = "5465489130"

$ws.Cells.Item(13, 3).Value = "asdjhf@hj.com"
$ws.Cells.Item(13, 4).Value = "ajshdfg"
